# "Generate Report for Handback"
# Updates the localization-status workbook to reflect a completed handback:
#  - Overview sheet: Status for both locales flips from "Ready for handoff"
#    to "Handed back: in sync with en-US".
#  - zh-cn / de-de sheets: each of the two tracked files now has a
#    "Latest Target File" (hyperlinked, like the Source File Name) and a
#    "Latest Handback File" populated, plus a real "Latest Handback DateTime".

$wb = $excel.ActiveWorkbook

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5969f903d479e5c85c4cb313497bac0eade0b00b/e2e/689a42b7-74b9-49ef-bc2f-6c2d1809b3e5.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5969f903d479e5c85c4cb313497bac0eade0b00b/e2e/ebd1250c-e908-4c4e-876f-884263b44b4f.md"

$nameA = "689a42b7-74b9-49ef-bc2f-6c2d1809b3e5.md"
$nameB = "ebd1250c-e908-4c4e-876f-884263b44b4f.md"

# ---------------------------------------------------------------------
# Overview sheet: status moves from "Ready for handoff" to the in-sync
# handback message for both locale columns / both files.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the two status columns to fit the new, longer text.
$wsOverview.Range("E1").ColumnWidth = 29.166666666666668
$wsOverview.Range("F1").ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet: populate Latest Target File / Latest Handback File /
# Latest Handback DateTime for both rows.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $urlA, "", "", $nameA)
$wsZhCn.Range("J2").Value = "689a42b7-74b9-49ef-bc2f-6c2d1809b3e5.d56561300d6ed0fde7e39eccdf7d5033262fc9dc.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-06 04:35:29"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $urlB, "", "", $nameB)
$wsZhCn.Range("J3").Value = "ebd1250c-e908-4c4e-876f-884263b44b4f.9ab32148eb48064a0f0e7be3a6b625af8d69cc8e.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-06 04:35:29"

$wsZhCn.Range("C1").ColumnWidth = 29.166666666666668
$wsZhCn.Range("I1").ColumnWidth = 39.17
$wsZhCn.Range("J1").ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: same shape of update, with its own handback timestamp.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $urlA, "", "", $nameA)
$wsDeDe.Range("J2").Value = "689a42b7-74b9-49ef-bc2f-6c2d1809b3e5.d56561300d6ed0fde7e39eccdf7d5033262fc9dc.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-06 04:35:48"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $urlB, "", "", $nameB)
$wsDeDe.Range("J3").Value = "ebd1250c-e908-4c4e-876f-884263b44b4f.9ab32148eb48064a0f0e7be3a6b625af8d69cc8e.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-06 04:35:48"

$wsDeDe.Range("C1").ColumnWidth = 29.166666666666668
$wsDeDe.Range("I1").ColumnWidth = 39.17
$wsDeDe.Range("J1").ColumnWidth = 39.17

Write-Output "Handback report generated"
